# ajout du dataset openfoodfact
#
# Slide 6, shape 2 ("Espace réservé du contenu 2") contains the run:
#   ") à 1 par exemple => nous verrons ceci  dans la partie «<NBSP>Exploration des Données avec les "
# which must become three runs:
#   ") à 1 par exemple => nous verrons ceci  dans la partie «<NBSP>Exploration "
#   "des Données<NBSP>» "
#   "avec les "
# i.e. a closing French guillemet ("<NBSP>» ") is inserted after "Données",
# splitting the original run into three runs at that point.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$nbsp = [char]0x00A0

# --- locate the target span ("des Données avec les ") inside the run ---
$full = $tr.Text
$idxDes = $full.IndexOf("des Donn")
$idxAvecLesEnd = $full.IndexOf("avec les ") + ("avec les ").Length

$start1 = $idxDes + 1
$len1 = $idxAvecLesEnd - $idxDes

# Rewrite that span, inserting the missing closing guillemet after "Données".
# Setting .Text on this sub-range also creates run boundaries at both ends,
# splitting the original single run into two runs.
$target = $tr.Characters($start1, $len1)
$target.Text = "des Données" + $nbsp + "» avec les "

# --- now split "des Données» avec les " into "des Données» " / "avec les " ---
$full2 = $tr.Text
$idxAvec = $full2.IndexOf("avec les ")

$start2 = $idxAvec + 1
$len2 = ("avec les ").Length

$target2 = $tr.Characters($start2, $len2)
$target2.Text = "avec les "
